# Refresh the "cryptos" price/volume snapshot (GitHub Actions cron update).
# Updates Price (D) / Volume(1h) (E) for each coin row, and also fixes the
# EthereumClassic/Stellar row ordering (rows 27-28 swap places along with
# their refreshed figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.401.74'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = '1.601.78'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''212.13'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = '''0.518'
$ws.Range('E6').Value = '  +6.47%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''26.54'
$ws.Range('E8').Value = '  +6.13%  '
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').Value = '1.830.59'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = '1.630.25'
$ws.Range('E14').Value = '  +4.07%  '
$ws.Range('D15').Value = '29.428.33'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('E16').Value = '  +3.68%  '
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('E18').Value = '  +2.94%  '
$ws.Range('D19').Value = '''240.21'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').Value = '''7.63'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').Value = '0.0₃0688'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '''3.99'
$ws.Range('E23').Value = '  +1.89%  '
$ws.Range('D24').Value = '''9.12'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').Value = '''2.09'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '''154.45'
$ws.Range('E26').Value = '  +2.87%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '''0.109'
$ws.Range('E27').Value = '  +4.82%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''15.25'
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('D29').Value = '''6.36'
$ws.Range('E29').Value = '  +2.24%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +2.34%  '
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  +3.58%  '
$ws.Range('D35').Value = '1.411.65'
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('E37').Value = '  +2.79%  '
$ws.Range('D38').Value = '''2.82'
$ws.Range('E38').Value = '  +5.46%  '
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').Value = '''0.535'
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').Value = '''0.0486'
$ws.Range('E43').Value = '  +5.48%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '''0.793'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').Value = '''52.56'
$ws.Range('E46').Value = '  +21.48%  '
$ws.Range('D47').Value = '''65.55'
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('D48').Value = '''5.25'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = '1.741.36'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').Value = '''0.853'
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').Value = '''86.56'
$ws.Range('E51').Value = '  +1.39%  '

# Cells whose new value looks numeric were written with a leading
# apostrophe so Excel keeps them as text (matching the source data's
# exact formatting, e.g. trailing zeros like "1.00"). Reset the style
# afterwards so no stray "quote prefix" cell style is left behind.
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
